$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy the formatting of column F (rows 1-56) onto the new column G,
# then set the header text and the numeric values.
$ws.Range("F1:F56").Copy()
$ws.Range("G1:G56").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# New column G width (target OOXML width="17"; ColumnWidth uses a
# slightly different pixel-rounded unit, so 16.14 round-trips to 17)
$ws.Columns.Item(7).ColumnWidth = 16.14
